# Add files via upload - "1.2 Test cases"
# Adds a new worksheet "1.2" (after Sheet1) containing a fresh set of
# Partner-Mgmt / IDA / ID-Repo test cases, and updates a couple of
# workbook-level view settings to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New worksheet, inserted right after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "1.2"

# ---------------------------------------------------------------------
# 2. Header row — copy Sheet1's header formatting (bold/fill/border
#    style) across, then stamp in the header captions.
# ---------------------------------------------------------------------
$ws1.Range("A1:G1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "TestCase_No "
$ws2.Range("B1").Value = "Module"
$ws2.Range("C1").Value = "Feature"
$ws2.Range("D1").Value = "Type"
$ws2.Range("E1").Value = "Scenario"
$ws2.Range("F1").Value = "Expected Result"
$ws2.Range("G1").Value = "Automated"

# ---------------------------------------------------------------------
# 3. Data rows (2-11)
# ---------------------------------------------------------------------
$rows = @(
    @("PM_Policy_01", "Partner Mgmt", "Policy", "Functional", "Create a policy with attribute staticTokenType ", "Policy Manager should be able to create the policy"),
    @("IDA_Auth_01", "IDA", "Auth", "Functional", "Create a policy with staticTokenType as 'Random' and perform auth", "Random token should be generated everytime"),
    @("IDA_Auth_02", "IDA", "Auth", "Functional", "Create a policy with staticTokenType as 'Partner' and perform auth", "Token should be generated using the Partner ID and UIN"),
    @("IDA_Auth_03", "IDA", "Auth", "Functional", "Create a policy with staticTokenType as 'Policy' and perform auth", "Token should be generated using the  Policy ID and UIN"),
    @("IDA_Auth_04", "IDA", "Auth", "Functional", "Create a policy with staticTokenType as 'Policy Group' and perform auth", "Token should be generated using the  Policy Group ID and UIN"),
    @("IDA_Auth_05", "IDA", "Auth", "Functional", "Verify auth request accepts consent token in the request", "Consent token should be passed as part of the auth request "),
    @("IDA_DB_01", "IDA", "DB", "Functional", "Verify uin_hash column in ida table has been changed to id_hash", "Column name should be changed to id_hash"),
    @("ID Repo_DB_01", "ID Repo", "DB", "Functional", "Verify uin_hash column in id repo table has been changed to id_hash", "Column name should be changed to id_hash"),
    @("ID Repo_Lock_01", "ID Repo", "Lock", "Functional", "Verify if the UIN is locked, the corresponding VID is also locked for auth", "The corresponding VID should also be  locked for auth"),
    @("ID Repo_Lock_02", "ID Repo", "Lock", "Functional", "Verify if the VID is locked, the corresponding UIN is also locked for auth", "The corresponding UIN should also be  locked for auth")
)

$r = 2
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4. Column widths (approximate the source workbook's custom widths)
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 16.02
$ws2.Columns.Item(2).ColumnWidth = 13.45
$ws2.Columns.Item(3).ColumnWidth = 15.02
$ws2.Columns.Item(4).ColumnWidth = 17.17
$ws2.Columns.Item(5).ColumnWidth = 63.74
$ws2.Columns.Item(6).ColumnWidth = 78.59
$ws2.Columns.Item(7).ColumnWidth = 11.59

# ---------------------------------------------------------------------
# 5. Data validations (drop-down lists) on the new sheet
# ---------------------------------------------------------------------
$ws2.Range("B2:B103").Validation.Add(3, 1, 1, '"Admin, IDA, ID Repo, Partner Mgmt, Pre Registration, Registration Client, Registration Processor, Resident Services"')
$ws2.Range("G2:G299").Validation.Add(3, 1, 1, '"Y,N"')
$ws2.Range("D2:D299").Validation.Add(3, 1, 1, '"Acceptance, Functional, Smoke, Security, Performance"')
$ws2.Range("B104:B299").Validation.Add(3, 1, 1, '"Admin, IDA, Partner Mgmt, Pre Registration, Registration Client, Registration Processor, Resident Services"')

# ---------------------------------------------------------------------
# 6. Selection / active-cell on the new sheet, and make it the active tab
# ---------------------------------------------------------------------
$ws2.Range("D14").Select()
$ws2.Activate()

Write-Output "done"
